$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new columns between existing A and B (new B, C for "A"/"AA" and "Good"/"Bad")
$ws.Range("B1:C2").Insert(-4161)

$ws.Range("B1").Value = "A"
$ws.Range("B2").Value = "AA"
$ws.Range("C1").Value = "Good"
$ws.Range("C2").Value = "Bad"

Write-Output "done"
